# PlayerPerformance_4507.xlsx update
# 1. Insert a new "Player Info" sheet before "ODI Batting" with player bio data.
# 2. Rename MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" and "ODI Bowling".
# 3. Replace the full scorecard URL values with the bare match-code number.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the "Player Info" sheet ahead of "ODI Batting"
# ---------------------------------------------------------------------------
$battingSheetForInsert = $wb.Worksheets.Item("ODI Batting")
$infoSheet = $wb.Worksheets.Add($battingSheetForInsert)
$infoSheet.Name = "Player Info"

$infoSheet.Range("A1").Value = "ID"
$infoSheet.Range("B1").Value = "NAME"
$infoSheet.Range("C1").Value = "BATTING_HAND"
$infoSheet.Range("D1").Value = "BOWL_STYLE"

$infoHeader = $infoSheet.Range("A1:D1")
$infoHeader.Font.Bold = $true
$infoHeader.HorizontalAlignment = -4108
$infoHeader.VerticalAlignment = -4160
$infoHeader.Borders.LineStyle = 1
$infoHeader.Borders.Weight = 2

$infoSheet.Range("A2").Value = "'4507"
$infoSheet.Range("A2").Style = "Normal"
$infoSheet.Range("B2").Value = "Yamin Ahmadzai"
$infoSheet.Range("C2").Value = "Right Handed"
$infoSheet.Range("D2").Value = "Right Arm Medium Fast"

$infoSheet.Range("A1").Select()

# ---------------------------------------------------------------------------
# Re-fetch the other two sheets by name now that the sheet collection has
# shifted (indices moved after the insert above) -- name lookups stay valid.
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# ---------------------------------------------------------------------------
# 2 & 3. ODI Batting: rename MATCH_CARD_LINK header -> MATCH_CODE, and
#        replace the scorecard URL in column D with the bare match code.
# ---------------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{
    2 = "3864"
    3 = "3917"
    4 = "3918"
    5 = "4379"
    6 = "4525"
    7 = "4528"
    8 = "4537"
    9 = "4671"
    10 = "4674"
}
foreach ($row in $battingCodes.Keys) {
    $cell = $battingSheet.Range("D$row")
    $cell.Value = "'" + $battingCodes[$row]
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# ODI Bowling: rename MATCH_CARD_LINK header -> MATCH_CODE, and replace the
#              scorecard URL in column B with the bare match code.
# ---------------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @{
    2 = "3864"
    3 = "3917"
    4 = "3918"
    5 = "4379"
    6 = "4525"
    7 = "4528"
    8 = "4537"
    9 = "4671"
}
foreach ($row in $bowlingCodes.Keys) {
    $cell = $bowlingSheet.Range("B$row")
    $cell.Value = "'" + $bowlingCodes[$row]
    $cell.Style = "Normal"
}
